$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("N").Delete()

$ws.Range("Q1").Value = "syst1_u"
$ws.Range("R1").Value = "syst2_u"
$ws.Range("W1").Value = "syst7_u"
$ws.Range("Y1").Value = "syst_tot"

$ws.Range("N27").Select()
